# "use blog in English and fix all CI/CD data"
# Translate the French header row of the NoraUi-blog sheet to English.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Title"
$ws.Range("C1").Value = "Text"
$ws.Range("D1").Value = "Author"
$ws.Range("F1").Value = "Result"

# Move the active cell selection from F6 to F2.
$ws.Range("F2").Select()
